# Append rows 8-13 to the "Artfynd" sheet (new species observation records).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 1).Value = 131106652
$ws.Cells.Item(8, 2).Value = 80308
$ws.Cells.Item(8, 4).Formula = "'LC"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = 229497
$ws.Cells.Item(8, 6).Formula = "'Korallblylav"
$ws.Cells.Item(8, 6).Style = "Normal"
$ws.Cells.Item(8, 7).Formula = "'Parmeliella triptophylla"
$ws.Cells.Item(8, 7).Style = "Normal"
$ws.Cells.Item(8, 8).Formula = "'(Ach.) Müll.Arg."
$ws.Cells.Item(8, 8).Style = "Normal"
$ws.Cells.Item(8, 9).Formula = "'"
$ws.Cells.Item(8, 9).Style = "Normal"
$ws.Cells.Item(8, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(8, 16).Style = "Normal"
$ws.Cells.Item(8, 17).Value = 601273
$ws.Cells.Item(8, 18).Value = 6959781
$ws.Cells.Item(8, 19).Value = 10
$ws.Cells.Item(8, 20).Formula = "'Västernorrland"
$ws.Cells.Item(8, 20).Style = "Normal"
$ws.Cells.Item(8, 21).Formula = "'Timrå"
$ws.Cells.Item(8, 21).Style = "Normal"
$ws.Cells.Item(8, 22).Formula = "'Medelpad"
$ws.Cells.Item(8, 22).Style = "Normal"
$ws.Cells.Item(8, 23).Formula = "'Ljustorp"
$ws.Cells.Item(8, 23).Style = "Normal"
$ws.Cells.Item(8, 24).Formula = "'2025_0524"
$ws.Cells.Item(8, 24).Style = "Normal"
$ws.Cells.Item(8, 25).Formula = "'2025-06-26"
$ws.Cells.Item(8, 25).Style = "Normal"
$ws.Cells.Item(8, 26).Formula = "'08:52"
$ws.Cells.Item(8, 26).Style = "Normal"
$ws.Cells.Item(8, 27).Formula = "'2025-06-26"
$ws.Cells.Item(8, 27).Style = "Normal"
$ws.Cells.Item(8, 28).Formula = "'08:52"
$ws.Cells.Item(8, 28).Style = "Normal"
$ws.Cells.Item(8, 29).Formula = "'aspstubbe"
$ws.Cells.Item(8, 29).Style = "Normal"
$ws.Cells.Item(8, 30).Value = $false
$ws.Cells.Item(8, 31).Value = $false
$ws.Cells.Item(8, 33).Value = $false
$ws.Cells.Item(8, 46).Formula = "'"
$ws.Cells.Item(8, 46).Style = "Normal"
$ws.Cells.Item(8, 49).Formula = "'David Isaksson"
$ws.Cells.Item(8, 49).Style = "Normal"
$ws.Cells.Item(8, 50).Formula = "'Måns Svensson"
$ws.Cells.Item(8, 50).Style = "Normal"
$ws.Cells.Item(8, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(8, 51).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 1).Value = 131106657
$ws.Cells.Item(9, 2).Value = 80348
$ws.Cells.Item(9, 4).Formula = "'NT"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = 6458
$ws.Cells.Item(9, 6).Formula = "'Lunglav"
$ws.Cells.Item(9, 6).Style = "Normal"
$ws.Cells.Item(9, 7).Formula = "'Lobaria pulmonaria"
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(9, 8).Formula = "'(L.) Hoffm."
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 9).Formula = "'1"
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 10).Formula = "'bålar"
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(9, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(9, 16).Style = "Normal"
$ws.Cells.Item(9, 17).Value = 601264
$ws.Cells.Item(9, 18).Value = 6959676
$ws.Cells.Item(9, 19).Value = 10
$ws.Cells.Item(9, 20).Formula = "'Västernorrland"
$ws.Cells.Item(9, 20).Style = "Normal"
$ws.Cells.Item(9, 21).Formula = "'Timrå"
$ws.Cells.Item(9, 21).Style = "Normal"
$ws.Cells.Item(9, 22).Formula = "'Medelpad"
$ws.Cells.Item(9, 22).Style = "Normal"
$ws.Cells.Item(9, 23).Formula = "'Ljustorp"
$ws.Cells.Item(9, 23).Style = "Normal"
$ws.Cells.Item(9, 24).Formula = "'2025_0519"
$ws.Cells.Item(9, 24).Style = "Normal"
$ws.Cells.Item(9, 25).Formula = "'2025-06-26"
$ws.Cells.Item(9, 25).Style = "Normal"
$ws.Cells.Item(9, 26).Formula = "'08:36"
$ws.Cells.Item(9, 26).Style = "Normal"
$ws.Cells.Item(9, 27).Formula = "'2025-06-26"
$ws.Cells.Item(9, 27).Style = "Normal"
$ws.Cells.Item(9, 28).Formula = "'08:36"
$ws.Cells.Item(9, 28).Style = "Normal"
$ws.Cells.Item(9, 29).Formula = "'gammal asp"
$ws.Cells.Item(9, 29).Style = "Normal"
$ws.Cells.Item(9, 30).Value = $false
$ws.Cells.Item(9, 31).Value = $false
$ws.Cells.Item(9, 33).Value = $false
$ws.Cells.Item(9, 46).Formula = "'"
$ws.Cells.Item(9, 46).Style = "Normal"
$ws.Cells.Item(9, 49).Formula = "'David Isaksson"
$ws.Cells.Item(9, 49).Style = "Normal"
$ws.Cells.Item(9, 50).Formula = "'Måns Svensson"
$ws.Cells.Item(9, 50).Style = "Normal"
$ws.Cells.Item(9, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(9, 51).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 1).Value = 131106656
$ws.Cells.Item(10, 2).Value = 80221
$ws.Cells.Item(10, 4).Formula = "'VU"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = 392
$ws.Cells.Item(10, 6).Formula = "'Aspgelélav"
$ws.Cells.Item(10, 6).Style = "Normal"
$ws.Cells.Item(10, 7).Formula = "'Collema subnigrescens"
$ws.Cells.Item(10, 7).Style = "Normal"
$ws.Cells.Item(10, 8).Formula = "'Degel."
$ws.Cells.Item(10, 8).Style = "Normal"
$ws.Cells.Item(10, 9).Formula = "'10"
$ws.Cells.Item(10, 9).Style = "Normal"
$ws.Cells.Item(10, 10).Formula = "'bålar"
$ws.Cells.Item(10, 10).Style = "Normal"
$ws.Cells.Item(10, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(10, 16).Style = "Normal"
$ws.Cells.Item(10, 17).Value = 601270
$ws.Cells.Item(10, 18).Value = 6959748
$ws.Cells.Item(10, 19).Value = 10
$ws.Cells.Item(10, 20).Formula = "'Västernorrland"
$ws.Cells.Item(10, 20).Style = "Normal"
$ws.Cells.Item(10, 21).Formula = "'Timrå"
$ws.Cells.Item(10, 21).Style = "Normal"
$ws.Cells.Item(10, 22).Formula = "'Medelpad"
$ws.Cells.Item(10, 22).Style = "Normal"
$ws.Cells.Item(10, 23).Formula = "'Ljustorp"
$ws.Cells.Item(10, 23).Style = "Normal"
$ws.Cells.Item(10, 24).Formula = "'2025_0520"
$ws.Cells.Item(10, 24).Style = "Normal"
$ws.Cells.Item(10, 25).Formula = "'2025-06-26"
$ws.Cells.Item(10, 25).Style = "Normal"
$ws.Cells.Item(10, 26).Formula = "'08:43"
$ws.Cells.Item(10, 26).Style = "Normal"
$ws.Cells.Item(10, 27).Formula = "'2025-06-26"
$ws.Cells.Item(10, 27).Style = "Normal"
$ws.Cells.Item(10, 28).Formula = "'08:43"
$ws.Cells.Item(10, 28).Style = "Normal"
$ws.Cells.Item(10, 29).Formula = "'asphögstubbe"
$ws.Cells.Item(10, 29).Style = "Normal"
$ws.Cells.Item(10, 30).Value = $false
$ws.Cells.Item(10, 31).Value = $false
$ws.Cells.Item(10, 33).Value = $false
$ws.Cells.Item(10, 46).Formula = "'"
$ws.Cells.Item(10, 46).Style = "Normal"
$ws.Cells.Item(10, 49).Formula = "'David Isaksson"
$ws.Cells.Item(10, 49).Style = "Normal"
$ws.Cells.Item(10, 50).Formula = "'Måns Svensson"
$ws.Cells.Item(10, 50).Style = "Normal"
$ws.Cells.Item(10, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(10, 51).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 1).Value = 131106655
$ws.Cells.Item(11, 2).Value = 80348
$ws.Cells.Item(11, 4).Formula = "'NT"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = 6458
$ws.Cells.Item(11, 6).Formula = "'Lunglav"
$ws.Cells.Item(11, 6).Style = "Normal"
$ws.Cells.Item(11, 7).Formula = "'Lobaria pulmonaria"
$ws.Cells.Item(11, 7).Style = "Normal"
$ws.Cells.Item(11, 8).Formula = "'(L.) Hoffm."
$ws.Cells.Item(11, 8).Style = "Normal"
$ws.Cells.Item(11, 9).Formula = "'1"
$ws.Cells.Item(11, 9).Style = "Normal"
$ws.Cells.Item(11, 10).Formula = "'bålar"
$ws.Cells.Item(11, 10).Style = "Normal"
$ws.Cells.Item(11, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(11, 16).Style = "Normal"
$ws.Cells.Item(11, 17).Value = 601279
$ws.Cells.Item(11, 18).Value = 6959779
$ws.Cells.Item(11, 19).Value = 10
$ws.Cells.Item(11, 20).Formula = "'Västernorrland"
$ws.Cells.Item(11, 20).Style = "Normal"
$ws.Cells.Item(11, 21).Formula = "'Timrå"
$ws.Cells.Item(11, 21).Style = "Normal"
$ws.Cells.Item(11, 22).Formula = "'Medelpad"
$ws.Cells.Item(11, 22).Style = "Normal"
$ws.Cells.Item(11, 23).Formula = "'Ljustorp"
$ws.Cells.Item(11, 23).Style = "Normal"
$ws.Cells.Item(11, 24).Formula = "'2025_0521"
$ws.Cells.Item(11, 24).Style = "Normal"
$ws.Cells.Item(11, 25).Formula = "'2025-06-26"
$ws.Cells.Item(11, 25).Style = "Normal"
$ws.Cells.Item(11, 26).Formula = "'08:44"
$ws.Cells.Item(11, 26).Style = "Normal"
$ws.Cells.Item(11, 27).Formula = "'2025-06-26"
$ws.Cells.Item(11, 27).Style = "Normal"
$ws.Cells.Item(11, 28).Formula = "'08:44"
$ws.Cells.Item(11, 28).Style = "Normal"
$ws.Cells.Item(11, 30).Value = $false
$ws.Cells.Item(11, 31).Value = $false
$ws.Cells.Item(11, 33).Value = $false
$ws.Cells.Item(11, 46).Formula = "'"
$ws.Cells.Item(11, 46).Style = "Normal"
$ws.Cells.Item(11, 49).Formula = "'David Isaksson"
$ws.Cells.Item(11, 49).Style = "Normal"
$ws.Cells.Item(11, 50).Formula = "'Samuel Koont"
$ws.Cells.Item(11, 50).Style = "Normal"
$ws.Cells.Item(11, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(11, 51).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 1).Value = 131106658
$ws.Cells.Item(12, 2).Value = 80221
$ws.Cells.Item(12, 4).Formula = "'VU"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = 392
$ws.Cells.Item(12, 6).Formula = "'Aspgelélav"
$ws.Cells.Item(12, 6).Style = "Normal"
$ws.Cells.Item(12, 7).Formula = "'Collema subnigrescens"
$ws.Cells.Item(12, 7).Style = "Normal"
$ws.Cells.Item(12, 8).Formula = "'Degel."
$ws.Cells.Item(12, 8).Style = "Normal"
$ws.Cells.Item(12, 9).Formula = "'1"
$ws.Cells.Item(12, 9).Style = "Normal"
$ws.Cells.Item(12, 10).Formula = "'bålar"
$ws.Cells.Item(12, 10).Style = "Normal"
$ws.Cells.Item(12, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(12, 16).Style = "Normal"
$ws.Cells.Item(12, 17).Value = 601282
$ws.Cells.Item(12, 18).Value = 6959785
$ws.Cells.Item(12, 19).Value = 10
$ws.Cells.Item(12, 20).Formula = "'Västernorrland"
$ws.Cells.Item(12, 20).Style = "Normal"
$ws.Cells.Item(12, 21).Formula = "'Timrå"
$ws.Cells.Item(12, 21).Style = "Normal"
$ws.Cells.Item(12, 22).Formula = "'Medelpad"
$ws.Cells.Item(12, 22).Style = "Normal"
$ws.Cells.Item(12, 23).Formula = "'Ljustorp"
$ws.Cells.Item(12, 23).Style = "Normal"
$ws.Cells.Item(12, 24).Formula = "'2025_0518"
$ws.Cells.Item(12, 24).Style = "Normal"
$ws.Cells.Item(12, 25).Formula = "'2025-06-26"
$ws.Cells.Item(12, 25).Style = "Normal"
$ws.Cells.Item(12, 26).Formula = "'08:35"
$ws.Cells.Item(12, 26).Style = "Normal"
$ws.Cells.Item(12, 27).Formula = "'2025-06-26"
$ws.Cells.Item(12, 27).Style = "Normal"
$ws.Cells.Item(12, 28).Formula = "'08:35"
$ws.Cells.Item(12, 28).Style = "Normal"
$ws.Cells.Item(12, 30).Value = $false
$ws.Cells.Item(12, 31).Value = $false
$ws.Cells.Item(12, 33).Value = $false
$ws.Cells.Item(12, 46).Formula = "'"
$ws.Cells.Item(12, 46).Style = "Normal"
$ws.Cells.Item(12, 49).Formula = "'David Isaksson"
$ws.Cells.Item(12, 49).Style = "Normal"
$ws.Cells.Item(12, 50).Formula = "'Samuel Koont"
$ws.Cells.Item(12, 50).Style = "Normal"
$ws.Cells.Item(12, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(12, 51).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 1).Value = 131106662
$ws.Cells.Item(13, 2).Value = 57881
$ws.Cells.Item(13, 4).Formula = "'NT"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = 100049
$ws.Cells.Item(13, 6).Formula = "'Spillkråka"
$ws.Cells.Item(13, 6).Style = "Normal"
$ws.Cells.Item(13, 7).Formula = "'Dryocopus martius"
$ws.Cells.Item(13, 7).Style = "Normal"
$ws.Cells.Item(13, 8).Formula = "'(Linnaeus, 1758)"
$ws.Cells.Item(13, 8).Style = "Normal"
$ws.Cells.Item(13, 9).Formula = "'1"
$ws.Cells.Item(13, 9).Style = "Normal"
$ws.Cells.Item(13, 16).Formula = "'Lappberget, Mpd"
$ws.Cells.Item(13, 16).Style = "Normal"
$ws.Cells.Item(13, 17).Value = 601246
$ws.Cells.Item(13, 18).Value = 6959852
$ws.Cells.Item(13, 19).Value = 10
$ws.Cells.Item(13, 20).Formula = "'Västernorrland"
$ws.Cells.Item(13, 20).Style = "Normal"
$ws.Cells.Item(13, 21).Formula = "'Timrå"
$ws.Cells.Item(13, 21).Style = "Normal"
$ws.Cells.Item(13, 22).Formula = "'Medelpad"
$ws.Cells.Item(13, 22).Style = "Normal"
$ws.Cells.Item(13, 23).Formula = "'Ljustorp"
$ws.Cells.Item(13, 23).Style = "Normal"
$ws.Cells.Item(13, 24).Formula = "'2025_0514"
$ws.Cells.Item(13, 24).Style = "Normal"
$ws.Cells.Item(13, 25).Formula = "'2025-06-26"
$ws.Cells.Item(13, 25).Style = "Normal"
$ws.Cells.Item(13, 26).Formula = "'08:25"
$ws.Cells.Item(13, 26).Style = "Normal"
$ws.Cells.Item(13, 27).Formula = "'2025-06-26"
$ws.Cells.Item(13, 27).Style = "Normal"
$ws.Cells.Item(13, 28).Formula = "'08:25"
$ws.Cells.Item(13, 28).Style = "Normal"
$ws.Cells.Item(13, 29).Formula = "'Hackspår av spillkråka"
$ws.Cells.Item(13, 29).Style = "Normal"
$ws.Cells.Item(13, 30).Value = $false
$ws.Cells.Item(13, 31).Value = $false
$ws.Cells.Item(13, 33).Value = $false
$ws.Cells.Item(13, 46).Formula = "'"
$ws.Cells.Item(13, 46).Style = "Normal"
$ws.Cells.Item(13, 49).Formula = "'David Isaksson"
$ws.Cells.Item(13, 49).Style = "Normal"
$ws.Cells.Item(13, 50).Formula = "'Samuel Koont"
$ws.Cells.Item(13, 50).Style = "Normal"
$ws.Cells.Item(13, 51).Formula = "'Kustpaketet"
$ws.Cells.Item(13, 51).Style = "Normal"
